# ISS seminar 2023 documentation fix ("ispravljena greska u dokumentaciji")
#
# The bulk of the source diff wraps foreign/English technical terms
# (Unity, Dalen, MonoBehaviour, Update, Instantiate, Destroy,
# OnCollisionEnter, Invoke, Assets, store, MortarShell, MortarAction,
# TankAction, TankShooting, SpawnerScript, Transform, Rigidbody,
# BoxCollider, AudioSource, JetBrains, s.r.o, IntelliJ, Creating/and/
# Using/Scripts, Explore/the/Unity/Editor, Unity Technologies, ...)
# in <w:proofErr w:type="spellStart"/> / <w:proofErr w:type="spellEnd"/>
# markers. Those markers are purely Word's internal spell-checker
# bookkeeping - they carry no visible text/content and are not part of
# the exposed Word object model (there both in real Word VBA/COM and
# in this runtime), so they cannot be produced from automation code;
# the surrounding run text is byte-for-byte identical before/after.
#
# The one substantive, content-visible change bundled in the same
# commit is a typo correction in the MortarAction/Fire paragraph:
# "vansjke" -> "vanjske".

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "vansjke",  # FindText
    $true,      # MatchCase
    $true,      # MatchWholeWord
    $false,     # MatchWildcards
    $false,     # MatchSoundsLike
    $false,     # MatchAllWordForms
    $true,      # Forward
    1,          # Wrap (wdFindContinue)
    $false,     # Format
    "vanjske",  # ReplaceWith
    2           # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Expected text 'vansjke' was not found in the document."
}
